$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1:B14").ClearContents()
$ws.Range("B1:B14").Select()
